$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.476.19'
$ws.Range("E2").Value = '  -0.62%  '
$ws.Range("D3").Value = '2.070.28'
$ws.Range("E3").Value = '  -0.14%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '232.30'
$ws.Range("E5").Value = '  -0.30%  '
$ws.Range("E6").Value = '  +0.99%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '57.31'
$ws.Range("E8").Value = '  -1.84%  '
$ws.Range("E9").Value = '  -1.13%  '
$ws.Range("E10").Value = '  -0.60%  '
$ws.Range("E11").Value = '  +1.58%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '14.83'
$ws.Range("E12").Value = '  +0.69%  '
$ws.Range("D13").Value = '2.376.28'
$ws.Range("E13").Value = '  -0.19%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.86'
$ws.Range("E14").Value = '  -0.23%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.764'
$ws.Range("E15").Value = '  -1.13%  '
$ws.Range("E16").Value = '  -0.70%  '
$ws.Range("D17").Value = '2.075.77'
$ws.Range("E17").Value = '  -0.18%  '
$ws.Range("D18").Value = '37.397.71'
$ws.Range("E18").Value = '  -0.68%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '70.42'
$ws.Range("E19").Value = '  -0.84%  '
$ws.Range("E20").Value = '  -2.98%  '
$ws.Range("D21").Value = '0.0₃0828'
$ws.Range("E21").Value = '  -0.75%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '228.13'
$ws.Range("E22").Value = '  +0.05%  '
$ws.Range("E23").Value = '  +0.10%  '
$ws.Range("E24").Value = '  -0.53%  '
$ws.Range("E25").Value = '  -1.17%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.61'
$ws.Range("E26").Value = '  +6.65%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '170.00'
$ws.Range("E27").Value = '  -0.56%  '
$ws.Range("E28").Value = '  -3.39%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.44'
$ws.Range("E29").Value = '  -0.01%  '
$ws.Range("E30").Value = '  -0.97%  '
$ws.Range("E31").Value = '  +1.13%  '
$ws.Range("E32").Value = '  -1.21%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0633'
$ws.Range("E33").Value = '  +0.25%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.62'
$ws.Range("E34").Value = '  -0.42%  '
$ws.Range("E35").Value = '  -0.58%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.81'
$ws.Range("E36").Value = '  -0.47%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.31'
$ws.Range("E37").Value = '  -2.33%  '
$ws.Range("E38").Value = '  +0.02%  '
$ws.Range("E39").Value = '  -0.80%  '
$ws.Range("E40").Value = '  +7.06%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '99.67'
$ws.Range("E41").Value = '  -0.96%  '
$ws.Range("E42").Value = '  +0.87%  '
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.20'
$ws.Range("E43").Value = '  +3.93%  '
$ws.Range("B44").Value = 'Cronos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0950'
$ws.Range("E44").Value = '  -2.58%  '
$ws.Range("D45").Value = '1.472.10'
$ws.Range("E45").Value = '  +2.17%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '16.76'
$ws.Range("E46").Value = '  +1.39%  '
$ws.Range("E47").Value = '  -1.42%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.96'
$ws.Range("E48").Value = '  -5.63%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.20'
$ws.Range("E49").Value = '  -2.96%  '
$ws.Range("E50").Value = '  -2.23%  '
$ws.Range("D51").Value = '2.259.00'
$ws.Range("E51").Value = '  -0.31%  '
